$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.994209885597229
$ws.Range("B1").Value = 2.680562019348145
$ws.Range("C1").Value = 4.723860263824463
$ws.Range("D1").Value = 2.323959589004517
$ws.Range("E1").Value = 1.072211027145386
